$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") { continue }

    $parts = $text -split ', '

    if ($parts.Count -gt 1 -and $parts[0] -eq "System" -and ($text -notlike "*admin@admin.com*")) {
        $rest = $parts[1..($parts.Count - 1)]
        $newParts = $rest + @("System")
        $newText = $newParts -join ', '
        $cell.Value = $newText
    }
}
